$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default "Normal" style with no special number format)
# used to restore default formatting on cells that need a NumberFormat="@"
# nudge to stop Excel auto-converting numeric-looking price strings into numbers.
$refStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = "66.296.25"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "3.526.28"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.44"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "  +4.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.65"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "  -4.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = $refStyle
$ws.Range("E7").Value = "  +5.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.638"
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  +5.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.82"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000277"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.32"
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "4.085.41"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "3.511.89"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "66.257.59"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.32"
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "  +7.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "85.79"
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.97"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "  +9.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.11"
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.11"
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.60"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "643.01"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.59"
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.76"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.156"
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = "  +13.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.64"
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0807"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.260.70"
$ws.Range("E41").Value = "  +8.81%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.93"
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.33"
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.55"
$ws.Range("D45").Style = $refStyle
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0419"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.74"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.44"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.44"
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = "  -0.12%  "
